$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3101.5952
$ws.Range("I15").Value = 3101.5952
$ws.Range("K15").Value = 9304.785600000001
$ws.Range("M15").Value = -9135.785600000001

$ws.Range("H98").Value = 2028.5834
$ws.Range("I98").Value = 1167.2727
$ws.Range("J98").Value = 11503
$ws.Range("K98").Value = 1167.2727
$ws.Range("L98").Value = 11503
$ws.Range("M98").Value = 330.7273
$ws.Range("N98").Value = -14499

$ws.Range("H116").Value = 2412.3333
$ws.Range("I116").Value = 1838.125
$ws.Range("J116").Value = 7006
$ws.Range("K116").Value = 1838.125
$ws.Range("L116").Value = 7006
$ws.Range("M116").Value = 1603.875
$ws.Range("N116").Value = -13890

$ws.Range("H122").Value = 2028.5834
$ws.Range("I122").Value = 1167.2727
$ws.Range("J122").Value = 11503
$ws.Range("K122").Value = 3501.8181
$ws.Range("L122").Value = 34509
$ws.Range("M122").Value = -1051.8181
$ws.Range("N122").Value = -39409

$ws.Range("H125").Value = 2363.111
$ws.Range("I125").Value = 2046.4
$ws.Range("J125").Value = 2759
$ws.Range("K125").Value = 18417.6
$ws.Range("L125").Value = 24831
$ws.Range("M125").Value = -15957.6
$ws.Range("N125").Value = -29751

$ws.Range("H127").Value = 125001060
$ws.Range("I127").Value = 200000350
$ws.Range("J127").Value = 2233.3333
$ws.Range("K127").Value = 600001050
$ws.Range("L127").Value = 6699.999899999999
$ws.Range("M127").Value = -599996090
$ws.Range("N127").Value = -16619.9999

$ws.Range("H132").Value = 5191.914
$ws.Range("I132").Value = 2568.44
$ws.Range("J132").Value = 11750.6
$ws.Range("K132").Value = 7705.32
$ws.Range("L132").Value = 35251.8
$ws.Range("M132").Value = -5175.32
$ws.Range("N132").Value = -40311.8

$ws.Range("H137").Value = 975833.8
$ws.Range("I137").Value = 1133.9796
$ws.Range("J137").Value = 6945870
$ws.Range("K137").Value = 3401.9388
$ws.Range("L137").Value = 20837610
$ws.Range("M137").Value = -851.9387999999999
$ws.Range("N137").Value = -20842710

$ws.Range("H138").Value = 3244891.8
$ws.Range("I138").Value = 1018.8571
$ws.Range("J138").Value = 7449912
$ws.Range("K138").Value = 3056.5713
$ws.Range("L138").Value = 22349736
$ws.Range("M138").Value = 2083.4287
$ws.Range("N138").Value = -22360016

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3393.64
$ws.Range("I32").Value = 2753.7234
$ws.Range("J32").Value = 13419
$ws.Range("K32").Value = 2753.7234
$ws.Range("L32").Value = 13419
$ws.Range("M32").Value = -2466.7234
$ws.Range("N32").Value = -13993

$ws.Range("H61").Value = 1762.1111
$ws.Range("I61").Value = 1687.4286
$ws.Range("J61").Value = 1866.6666
$ws.Range("K61").Value = 1687.4286
$ws.Range("L61").Value = 1866.6666
$ws.Range("M61").Value = -1475.4286
$ws.Range("N61").Value = -2290.6666

$ws.Range("H74").Value = 22035.547
$ws.Range("I74").Value = 28024.73
$ws.Range("J74").Value = 8185.5625
$ws.Range("K74").Value = 28024.73
$ws.Range("L74").Value = 8185.5625
$ws.Range("M74").Value = -27150.73
$ws.Range("N74").Value = -9933.5625

$ws.Range("H77").Value = 22035.547
$ws.Range("I77").Value = 28024.73
$ws.Range("J77").Value = 8185.5625
$ws.Range("K77").Value = 140123.65
$ws.Range("L77").Value = 40927.8125
$ws.Range("M77").Value = -135755.65
$ws.Range("N77").Value = -49663.8125

$ws.Range("H110").Value = 1840.9524
$ws.Range("I110").Value = 1925.5555
$ws.Range("J110").Value = 1333.3334
$ws.Range("K110").Value = 1925.5555
$ws.Range("L110").Value = 1333.3334
$ws.Range("M110").Value = 119.4445000000001
$ws.Range("N110").Value = -5423.3334

$ws.Range("H136").Value = 1762.1111
$ws.Range("I136").Value = 1687.4286
$ws.Range("J136").Value = 1866.6666
$ws.Range("K136").Value = 5062.2858
$ws.Range("L136").Value = 5599.9998
$ws.Range("M136").Value = -2512.2858
$ws.Range("N136").Value = -10699.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1607.1389
$ws.Range("I20").Value = 1273.8096
$ws.Range("J20").Value = 2073.8
$ws.Range("K20").Value = 1273.8096
$ws.Range("L20").Value = 2073.8
$ws.Range("M20").Value = -1026.8096
$ws.Range("N20").Value = -2567.8

$ws.Range("H134").Value = 745378.6
$ws.Range("I134").Value = 1338158.1
$ws.Range("J134").Value = 4404.25
$ws.Range("K134").Value = 4014474.3
$ws.Range("L134").Value = 13212.75
$ws.Range("M134").Value = -4011939.3
$ws.Range("N134").Value = -18282.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23811112
$ws.Range("I31").Value = 1491.5294
$ws.Range("J31").Value = 125002000
$ws.Range("K31").Value = 1491.5294
$ws.Range("L31").Value = 125002000
$ws.Range("M31").Value = -1196.5294
$ws.Range("N31").Value = -125002590

$ws.Range("H34").Value = 23811112
$ws.Range("I34").Value = 1491.5294
$ws.Range("J34").Value = 125002000
$ws.Range("K34").Value = 1491.5294
$ws.Range("L34").Value = 125002000
$ws.Range("M34").Value = -1289.5294
$ws.Range("N34").Value = -125002404

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 851.56525
$ws.Range("I5").Value = 528.2857
$ws.Range("J5").Value = 1354.4445
$ws.Range("K5").Value = 1584.8571
$ws.Range("L5").Value = 4063.3335
$ws.Range("M5").Value = -1472.8571
$ws.Range("N5").Value = -4287.333500000001

$ws.Range("H56").Value = 4156.364
$ws.Range("I56").Value = 4156.364
$ws.Range("K56").Value = 4156.364
$ws.Range("M56").Value = -3626.364

$ws.Range("H113").Value = 583168.9
$ws.Range("I113").Value = 977935.0600000001
$ws.Range("J113").Value = 418.7619
$ws.Range("K113").Value = 2933805.18
$ws.Range("L113").Value = 1256.2857
$ws.Range("M113").Value = -2931635.18
$ws.Range("N113").Value = -5596.2857

$ws.Range("H122").Value = 58634.42
$ws.Range("I122").Value = 73870.336
$ws.Range("K122").Value = 664833.024
$ws.Range("M122").Value = -662383.024

$ws.Range("H135").Value = 851.56525
$ws.Range("I135").Value = 528.2857
$ws.Range("J135").Value = 1354.4445
$ws.Range("K135").Value = 4754.571300000001
$ws.Range("L135").Value = 12190.0005
$ws.Range("M135").Value = -2219.571300000001
$ws.Range("N135").Value = -17260.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2043402.5
$ws.Range("I132").Value = 2531.36
$ws.Range("K132").Value = 7594.08
$ws.Range("M132").Value = -5064.08

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3649.75
$ws.Range("I7").Value = 2299.5
$ws.Range("K7").Value = 2299.5
$ws.Range("M7").Value = -2187.5

$ws.Range("H40").Value = 2573.5454
$ws.Range("I40").Value = 2325.5
$ws.Range("J40").Value = 3235
$ws.Range("K40").Value = 2325.5
$ws.Range("L40").Value = 3235
$ws.Range("M40").Value = -2189.5
$ws.Range("N40").Value = -3507

$ws.Range("H126").Value = 3649.75
$ws.Range("I126").Value = 2299.5
$ws.Range("K126").Value = 6898.5
$ws.Range("M126").Value = -4428.5

$ws.Range("H132").Value = 3316.5078
$ws.Range("I132").Value = 2903.2979
$ws.Range("J132").Value = 4395.4443
$ws.Range("K132").Value = 8709.893700000001
$ws.Range("L132").Value = 13186.3329
$ws.Range("M132").Value = -6179.893700000001
$ws.Range("N132").Value = -18246.3329

$ws.Range("H136").Value = 2851.1667
$ws.Range("I136").Value = 1534
$ws.Range("J136").Value = 4168.3335
$ws.Range("K136").Value = 4602
$ws.Range("L136").Value = 12505.0005
$ws.Range("M136").Value = -2052
$ws.Range("N136").Value = -17605.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1728.5714
$ws.Range("I14").Value = 1050
$ws.Range("K14").Value = 1050
$ws.Range("M14").Value = -882

$ws.Range("H126").Value = 41667732
$ws.Range("I126").Value = 841.0769
$ws.Range("J126").Value = 90910424
$ws.Range("K126").Value = 2523.2307
$ws.Range("L126").Value = 272731272
$ws.Range("M126").Value = -53.23070000000007
$ws.Range("N126").Value = -272736212

$ws.Range("H132").Value = 2031.5231
$ws.Range("I132").Value = 1862.7872
$ws.Range("J132").Value = 2472.111
$ws.Range("K132").Value = 5588.3616
$ws.Range("L132").Value = 7416.333
$ws.Range("M132").Value = -3058.3616
$ws.Range("N132").Value = -12476.333

$ws.Range("H136").Value = 6120.8667
$ws.Range("I136").Value = 7338.625
$ws.Range("J136").Value = 4729.143
$ws.Range("K136").Value = 22015.875
$ws.Range("L136").Value = 14187.429
$ws.Range("M136").Value = -19465.875
$ws.Range("N136").Value = -19287.429
